# Consolidation_Evaluation.xlsx - reset the score inputs and add the
# "tasks" legend block (columns F:H) next to the evaluation grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Clear all of the previously-entered scores (columns B:D, rows 2-22).
# The E-column averages are formulas and recalculate to 0 automatically.
$ws.Range("B2:D22").ClearContents() | Out-Null

# --- New "tasks" legend block (F2, G2:G4, H2:H4) - written in the same
# order the original authoring tool produced them so the shared-string
# table comes out in the same sequence.
$ws.Range("F2").Value = " "
$ws.Range("G2").Value = "Tâche 1 : "
$ws.Range("G4").Value = "Tâche 3 :"
$ws.Range("G3").Value = "Tâche 2 : "
$ws.Range("H2").Value = "Réserver un cours"
$ws.Range("H3").Value = "Rechercher de nouveaux tournois"
$ws.Range("H4").Value = "Regarder des photos des évènements passés"

# --- Column sizing / outline grouping for the new layout.
$ws.Columns("C").ColumnWidth = 10.833333333333334
$ws.Columns("D").ColumnWidth = 10.833333333333334
$ws.Columns("D").Hidden = $true
$ws.Range("D1").EntireColumn.OutlineLevel = 1
$ws.Columns("G").ColumnWidth = 8.0
$ws.Columns("H").ColumnWidth = 40.5

# --- Move the active selection to B2 (first score cell) as in the saved file.
$ws.Range("B2").Select() | Out-Null

Write-Output "done"
